# Remove 0344, 3630, and FI for 0360
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "International Tracked Packet" / service code 3630 row entirely.
# (Sheet row 57: Service name "International Tracked Packet", Service code/Request code "3630".)
$ws.Rows(57).Delete()

# Re-establish the AutoFilter over the now-smaller data range (was A1:P63,
# now one row shorter at A1:P62) since deleting the row does not
# automatically shrink the existing filter range.
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:P62").AutoFilter()

# The hidden _FilterDatabase defined name also needs to be pointed at the
# new, smaller range.
$fdb = $wb.Names.Item(1)
$fdb.RefersTo = "='Booking & SG API'!`$A`$1:`$P`$62"

# Remove Finland (FI) from the sender-countries / destination list for
# service 0360 (Bring Pack, row 27): "SE, DK, FI" -> "SE, DK".
$ws.Range("O27").Value = "SE, DK"
$ws.Range("P27").Value = "SE, DK"
